# feat: add 2022-Q1 data
#
# - The existing "总计" sheet (sheetId 6) is renamed to "2022-Q1" and its
#   content is replaced by the per-fund holdings table for the new quarter.
# - A brand new "总计" sheet is appended after it, reusing the previous
#   "总计" rows (2021-Q4 .. 2020-Q4) with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Turn the old "总计" sheet into the new "2022-Q1" per-fund sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$fundHeader = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

$fundRows = @(
    @("001320", "工银瑞信丰盈回报灵活配置混合", "6.50", "93.76", "3.49", "0.2268", 9),
    @("000763", "工银新财富灵活配置混合", "2.96", "92.68", "3.91", "0.1157", 5),
    @("010617", "兴业消费精选混合A", "2.96", "68.95", "3.81", "0.1128", 9),
    @("010618", "兴业消费精选混合C", "1.47", "68.95", "3.81", "0.0560", 9),
    @("001692", "南方国策动力股票", "2.14", "94.09", "2.16", "0.0462", 7),
    @("011858", "安信消费升级一年持有期混合型发起式证券投资基金A", "1.10", "72.03", "2.89", "0.0318", 8),
    @("011500", "九泰量化新兴产业混合型证券投资基金", "0.70", "93.48", "2.51", "0.0176", 4),
    @("008353", "泰达宏利消费行业量化精选混合A", "0.52", "92.25", "1.65", "0.0086", 9),
    @("165524", "信诚中证智能家居指数（LOF）", "0.40", "93.89", "1.18", "0.0047", 8),
    @("011859", "安信消费升级一年持有期混合型发起式证券投资基金C", "0.14", "72.03", "2.89", "0.0040", 8),
    @("002330", "兴业聚宝灵活配置混合", "0.05", "80.30", "3.97", "0.0020", 5),
    @("008354", "泰达宏利消费行业量化精选混合C", "0.12", "92.25", "1.65", "0.0020", 9)
)

# Header row (B1:H1) copies the bold/centered/bordered look used on every
# other quarter sheet by reusing that sheet's own header style.
$headerStyleSrc = $wb.Worksheets.Item("2021-Q4").Cells.Item(1, 2)
for ($i = 0; $i -lt $fundHeader.Count; $i++) {
    $cell = $q1.Cells.Item(1, $i + 2)
    $headerStyleSrc.Copy($cell)
    $cell.Value = $fundHeader[$i]
}

# Columns B..G hold text (fund code / name / scale / position / ratio /
# value) even though most look numeric (e.g. "6.50", "001320" with a
# leading zero) - format as Text first so COM keeps them as strings
# instead of silently parsing them into doubles and dropping the
# formatting / leading zeros.
$q1.Range("B2:G13").NumberFormat = "@"

$colStyleSrc = $wb.Worksheets.Item("2021-Q4").Cells.Item(2, 1)
for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $aCell = $q1.Cells.Item($r, 1)
    $colStyleSrc.Copy($aCell)
    $aCell.Value = $i

    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2. Append a fresh "总计" sheet after "2022-Q1" with the summary table.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# Match the page-setup margins used by every other sheet in the workbook
# (0.75in sides, 1in top/bottom, 0.5in header/footer).
$total.PageSetup.LeftMargin = 0.75 * 72
$total.PageSetup.RightMargin = 0.75 * 72
$total.PageSetup.TopMargin = 1 * 72
$total.PageSetup.BottomMargin = 1 * 72
$total.PageSetup.HeaderMargin = 0.5 * 72
$total.PageSetup.FooterMargin = 0.5 * 72

$totalHeader = @("日期", "持有数量(只)", "持有市值(亿元)")
$totalRows = @(
    @("2022-Q1", 12, 0.63),
    @("2021-Q4", 31, 7.56),
    @("2021-Q3", 14, 1.33),
    @("2021-Q2", 29, 5.1),
    @("2021-Q1", 25, 8.109999999999999),
    @("2020-Q4", 14, 7.64)
)

for ($i = 0; $i -lt $totalHeader.Count; $i++) {
    $cell = $total.Cells.Item(1, $i + 2)
    $headerStyleSrc.Copy($cell)
    $cell.Value = $totalHeader[$i]
}

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    $aCell = $total.Cells.Item($r, 1)
    $colStyleSrc.Copy($aCell)
    $aCell.Value = $i

    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

Write-Output "done"
